$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the price/volume columns remain stored as text (matching the
# original inline-string cells) rather than being auto-converted to numbers.
$dataRange = $ws.Range("D2:E50")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "301.10"
$ws.Range("E2").Value = "0.45%"
$ws.Range("E3").Value = "0.93%"
$ws.Range("D4").Value = "5.074"
$ws.Range("E4").Value = "-1.10%"
$ws.Range("D5").Value = "0.07844"
$ws.Range("E5").Value = "-3.37%"
$ws.Range("D6").Value = "2.342"
$ws.Range("E6").Value = "-10.82%"
$ws.Range("D7").Value = "7.819"
$ws.Range("E7").Value = "-0.41%"
$ws.Range("D8").Value = "3.837"
$ws.Range("E8").Value = "-0.11%"
$ws.Range("D9").Value = "0.9179"
$ws.Range("E9").Value = "0.98%"
$ws.Range("D10").Value = "0.1758"
$ws.Range("E10").Value = "2.29%"
$ws.Range("D11").Value = "0.07592"
$ws.Range("E11").Value = "4.68%"
$ws.Range("D12").Value = "0.09263"
$ws.Range("E12").Value = "16.21%"
$ws.Range("D13").Value = "0.02989"
$ws.Range("E13").Value = "-1.20%"
$ws.Range("D14").Value = "0.1001"
$ws.Range("E14").Value = "0.44%"
$ws.Range("D15").Value = "0.001510"
$ws.Range("E15").Value = "0.71%"
$ws.Range("D16").Value = "0.005841"
$ws.Range("E16").Value = "-2.90%"
$ws.Range("D17").Value = "3.472"
$ws.Range("E17").Value = "-0.72%"
$ws.Range("D18").Value = "2.247"
$ws.Range("E18").Value = "-0.30%"
$ws.Range("E19").Value = "0.48%"
$ws.Range("E20").Value = "-0.91%"
$ws.Range("D21").Value = "4.053"
$ws.Range("E21").Value = "-11.68%"
$ws.Range("D22").Value = "0.1789"
$ws.Range("E22").Value = "11.64%"
$ws.Range("D23").Value = "0.04615"
$ws.Range("E23").Value = "0.96%"
$ws.Range("D24").Value = "0.001250"
$ws.Range("E24").Value = "-1.29%"
$ws.Range("D25").Value = "0.004475"
$ws.Range("E25").Value = "0.65%"
$ws.Range("E26").Value = "5.78%"
$ws.Range("E27").Value = "-1.50%"
$ws.Range("D39").Value = "0.01765"
$ws.Range("E39").Value = "-2.88%"
$ws.Range("D40").Value = "0.04781"
$ws.Range("E40").Value = "5.35%"
$ws.Range("D41").Value = "0.007208"
$ws.Range("E41").Value = "2.20%"
$ws.Range("D42").Value = "0.1360"
$ws.Range("E42").Value = "1.25%"
$ws.Range("E43").Value = "-2.37%"
$ws.Range("D44").Value = "0.01033"
$ws.Range("E44").Value = "-1.63%"
$ws.Range("D45").Value = "0.00006265"
$ws.Range("E45").Value = "-0.70%"
$ws.Range("E46").Value = "-0.13%"
$ws.Range("E47").Value = "24.66%"
$ws.Range("D48").Value = "0.7433"
$ws.Range("E48").Value = "-9.42%"
$ws.Range("D49").Value = "0.00002099"
$ws.Range("E49").Value = "-0.13%"
$ws.Range("D50").Value = "0.0001999"
$ws.Range("E50").Value = "-0.13%"

# Restore default (General) formatting/style so cells match the original look.
$dataRange.NumberFormat = "General"
$dataRange.Style = "Normal"

